$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 158, shifting existing rows 158..254 down to 159..255
$ws.Rows.Item(158).Insert()

# Populate the newly inserted row 158 with the new weekly record
$ws.Cells.Item(158, 1).Value = 6
$ws.Cells.Item(158, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(158, 3).Value = "Metropolitana"
$ws.Cells.Item(158, 4).Value = 44879
$ws.Cells.Item(158, 5).Value = 13
$ws.Cells.Item(158, 6).Value = 100112029
$ws.Cells.Item(158, 7).Value = "Orégano"
$ws.Cells.Item(158, 8).Value = "Sin especificar"
$ws.Cells.Item(158, 9).Value = "Primera"
$ws.Cells.Item(158, 10).Value = 41
$ws.Cells.Item(158, 11).Value = 16000
$ws.Cells.Item(158, 12).Value = 17000
$ws.Cells.Item(158, 13).Value = 16439
$ws.Cells.Item(158, 14).Value = "$/docena de atados"
$ws.Cells.Item(158, 15).Value = "Región Metropolitana"
$ws.Cells.Item(158, 16).Value = 5480
$ws.Cells.Item(158, 17).Value = 3
$ws.Cells.Item(158, 18).Value = "Hortaliza"
